$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.144.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.477.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.72%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -1.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.70"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.10%  "
$ws.Range("E10").Value = "  -1.19%  "
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.072.22"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("E14").Value = "  -2.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.480.05"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.124.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.15%  "
$ws.Range("E19").Value = "  -1.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "384.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.567"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.37%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.618.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000112"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("E29").Value = "  -0.49%  "
$ws.Range("E30").Value = "  -3.73%  "
$ws.Range("E31").Value = "  -3.58%  "
$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.155"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.96%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.43"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.505.94"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.21%  "
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "22.97"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.89%  "
$ws.Range("E37").Value = "  +1.15%  "
$ws.Range("E38").Value = "  -1.89%  "
$ws.Range("E39").Value = "  -3.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "161.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0778"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.798"
$ws.Range("D42").Style = "Normal"
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.52"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.81%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.33"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.40%  "
$ws.Range("E46").Value = "  -2.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.02%  "
$ws.Range("E48").Value = "  -3.14%  "
$ws.Range("E49").Value = "  -0.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.900"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.342.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.76%  "
